$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
$ws.Cells.Item(32, 8).Value = 967.3333
$ws.Cells.Item(32, 9).Value = 235.33333
$ws.Cells.Item(32, 11).Value = 235.33333
$ws.Cells.Item(32, 13).Value = 90.66667000000001
$ws.Cells.Item(43, 8).Value = 9236.75
$ws.Cells.Item(43, 10).Value = 10649
$ws.Cells.Item(43, 12).Value = 10649
$ws.Cells.Item(43, 14).Value = -10787
$ws.Cells.Item(58, 8).Value = 1463.3572
$ws.Cells.Item(58, 9).Value = 248.7
$ws.Cells.Item(58, 10).Value = 4500
$ws.Cells.Item(58, 11).Value = 746.0999999999999
$ws.Cells.Item(58, 12).Value = 13500
$ws.Cells.Item(58, 13).Value = -596.0999999999999
$ws.Cells.Item(58, 14).Value = -13800
$ws.Cells.Item(61, 8).Value = 3627.5
$ws.Cells.Item(61, 9).Value = 3353
$ws.Cells.Item(61, 10).Value = 5000
$ws.Cells.Item(61, 11).Value = 10059
$ws.Cells.Item(61, 12).Value = 15000
$ws.Cells.Item(61, 13).Value = -9887
$ws.Cells.Item(61, 14).Value = -15344
$ws.Cells.Item(70, 8).Value = 3374.25
$ws.Cells.Item(70, 10).Value = 3998.5
$ws.Cells.Item(70, 12).Value = 11995.5
$ws.Cells.Item(70, 14).Value = -12535.5
$ws.Cells.Item(73, 8).Value = 3374.25
$ws.Cells.Item(73, 10).Value = 3998.5
$ws.Cells.Item(73, 12).Value = 11995.5
$ws.Cells.Item(73, 14).Value = -13867.5
$ws.Cells.Item(106, 8).Value = 5237
$ws.Cells.Item(106, 9).Value = 5384.4
$ws.Cells.Item(106, 10).Value = 4500
$ws.Cells.Item(106, 11).Value = 5384.4
$ws.Cells.Item(106, 12).Value = 4500
$ws.Cells.Item(106, 13).Value = -4753.4
$ws.Cells.Item(106, 14).Value = -5762
$ws.Cells.Item(112, 8).Value = 6282.9
$ws.Cells.Item(112, 10).Value = 3319.25
$ws.Cells.Item(112, 12).Value = 9957.75
$ws.Cells.Item(112, 14).Value = -12173.75
$ws.Cells.Item(115, 8).Value = 550.75
$ws.Cells.Item(115, 10).Value = 499
$ws.Cells.Item(115, 12).Value = 1497
$ws.Cells.Item(115, 14).Value = -4631
$ws.Cells.Item(124, 8).Value = 250000
$ws.Cells.Item(124, 10).Value = 250000
$ws.Cells.Item(124, 12).Value = 250000
$ws.Cells.Item(124, 14).Value = -259820
$ws.Cells.Item(138, 8).Value = 3557.5676
$ws.Cells.Item(138, 9).Value = 936.4545
$ws.Cells.Item(138, 10).Value = 4666.5
$ws.Cells.Item(138, 11).Value = 2809.3635
$ws.Cells.Item(138, 12).Value = 13999.5
$ws.Cells.Item(138, 13).Value = 2330.6365
$ws.Cells.Item(138, 14).Value = -24279.5

$ws = $wb.Worksheets("ARM")
$ws.Cells.Item(2, 8).Value = 1458.6
$ws.Cells.Item(2, 9).Value = 1334.4
$ws.Cells.Item(2, 10).Value = 1582.8
$ws.Cells.Item(2, 11).Value = 1334.4
$ws.Cells.Item(2, 12).Value = 1582.8
$ws.Cells.Item(2, 13).Value = -1221.4
$ws.Cells.Item(2, 14).Value = -1808.8
$ws.Cells.Item(61, 8).Value = 1096452.8
$ws.Cells.Item(61, 9).Value = 3405.7058
$ws.Cells.Item(61, 10).Value = 2644936
$ws.Cells.Item(61, 11).Value = 3405.7058
$ws.Cells.Item(61, 12).Value = 2644936
$ws.Cells.Item(61, 13).Value = -3193.7058
$ws.Cells.Item(61, 14).Value = -2645360
$ws.Cells.Item(63, 8).Value = 2239.1
$ws.Cells.Item(63, 9).Value = 2341.5715
$ws.Cells.Item(63, 11).Value = 2341.5715
$ws.Cells.Item(63, 13).Value = -1655.5715
$ws.Cells.Item(66, 8).Value = 2239.1
$ws.Cells.Item(66, 9).Value = 2341.5715
$ws.Cells.Item(66, 11).Value = 11707.8575
$ws.Cells.Item(66, 13).Value = -8275.8575
$ws.Cells.Item(116, 8).Value = 1458.6
$ws.Cells.Item(116, 9).Value = 1334.4
$ws.Cells.Item(116, 10).Value = 1582.8
$ws.Cells.Item(116, 11).Value = 1334.4
$ws.Cells.Item(116, 12).Value = 1582.8
$ws.Cells.Item(116, 13).Value = 959.5999999999999
$ws.Cells.Item(116, 14).Value = -6170.8
$ws.Cells.Item(132, 8).Value = 5095723
$ws.Cells.Item(132, 9).Value = 1395.6818
$ws.Cells.Item(132, 11).Value = 4187.0454
$ws.Cells.Item(132, 13).Value = -1657.0454
$ws.Cells.Item(133, 8).Value = 69995
$ws.Cells.Item(133, 10).Value = 69995
$ws.Cells.Item(133, 12).Value = 69995
$ws.Cells.Item(133, 14).Value = -75055
$ws.Cells.Item(136, 8).Value = 1096452.8
$ws.Cells.Item(136, 9).Value = 3405.7058
$ws.Cells.Item(136, 10).Value = 2644936
$ws.Cells.Item(136, 11).Value = 10217.1174
$ws.Cells.Item(136, 12).Value = 7934808
$ws.Cells.Item(136, 13).Value = -7667.117400000001
$ws.Cells.Item(136, 14).Value = -7939908

$ws = $wb.Worksheets("BSM")
$ws.Cells.Item(3, 8).Value = 1458.6
$ws.Cells.Item(3, 9).Value = 1334.4
$ws.Cells.Item(3, 10).Value = 1582.8
$ws.Cells.Item(3, 11).Value = 1334.4
$ws.Cells.Item(3, 12).Value = 1582.8
$ws.Cells.Item(3, 13).Value = -1220.4
$ws.Cells.Item(3, 14).Value = -1810.8
$ws.Cells.Item(94, 8).Value = 1394.2778
$ws.Cells.Item(94, 10).Value = 1661.9166
$ws.Cells.Item(94, 12).Value = 1661.9166
$ws.Cells.Item(94, 14).Value = -2563.9166

$ws = $wb.Worksheets("CRP")
$ws.Cells.Item(31, 8).Value = 7223.515
$ws.Cells.Item(31, 9).Value = 1024.3572
$ws.Cells.Item(31, 10).Value = 41938.8
$ws.Cells.Item(31, 11).Value = 1024.3572
$ws.Cells.Item(31, 12).Value = 41938.8
$ws.Cells.Item(31, 13).Value = -729.3571999999999
$ws.Cells.Item(31, 14).Value = -42528.8
$ws.Cells.Item(34, 8).Value = 7223.515
$ws.Cells.Item(34, 9).Value = 1024.3572
$ws.Cells.Item(34, 10).Value = 41938.8
$ws.Cells.Item(34, 11).Value = 1024.3572
$ws.Cells.Item(34, 12).Value = 41938.8
$ws.Cells.Item(34, 13).Value = -822.3571999999999
$ws.Cells.Item(34, 14).Value = -42342.8
$ws.Cells.Item(75, 8).Value = 10840
$ws.Cells.Item(78, 8).Value = 10840
$ws.Cells.Item(86, 8).Value = 29997
$ws.Cells.Item(86, 9).Value = 49995
$ws.Cells.Item(86, 10).Value = 9999
$ws.Cells.Item(86, 11).Value = 49995
$ws.Cells.Item(86, 12).Value = 9999
$ws.Cells.Item(86, 13).Value = -48872
$ws.Cells.Item(86, 14).Value = -12245
$ws.Cells.Item(89, 8).Value = 29997
$ws.Cells.Item(89, 9).Value = 49995
$ws.Cells.Item(89, 10).Value = 9999
$ws.Cells.Item(89, 11).Value = 249975
$ws.Cells.Item(89, 12).Value = 49995
$ws.Cells.Item(89, 13).Value = -244359
$ws.Cells.Item(89, 14).Value = -61227
$ws.Cells.Item(107, 8).Value = 1721
$ws.Cells.Item(107, 9).Value = 1444.1364
$ws.Cells.Item(107, 10).Value = 2591.1428
$ws.Cells.Item(107, 11).Value = 1444.1364
$ws.Cells.Item(107, 12).Value = 2591.1428
$ws.Cells.Item(107, 13).Value = 475.8635999999999
$ws.Cells.Item(107, 14).Value = -6431.1428
$ws.Cells.Item(132, 8).Value = 29414034
$ws.Cells.Item(132, 9).Value = 1942.8572
$ws.Cells.Item(132, 10).Value = 147062400
$ws.Cells.Item(132, 11).Value = 5828.571599999999
$ws.Cells.Item(132, 12).Value = 441187200
$ws.Cells.Item(132, 13).Value = -3298.571599999999
$ws.Cells.Item(132, 14).Value = -441192260

$ws = $wb.Worksheets("CUL")
$ws.Cells.Item(4, 8).Value = 2507959
$ws.Cells.Item(4, 9).Value = 2809550.8
$ws.Cells.Item(4, 11).Value = 8428652.399999999
$ws.Cells.Item(4, 13).Value = -8428540.399999999
$ws.Cells.Item(7, 8).Value = 6842.625
$ws.Cells.Item(7, 9).Value = 120.5
$ws.Cells.Item(7, 10).Value = 9083.333
$ws.Cells.Item(7, 11).Value = 361.5
$ws.Cells.Item(7, 12).Value = 27249.999
$ws.Cells.Item(7, 13).Value = -249.5
$ws.Cells.Item(7, 14).Value = -27473.999
$ws.Cells.Item(68, 8).Value = 1033.375
$ws.Cells.Item(68, 9).Value = 899
$ws.Cells.Item(68, 10).Value = 1052.5714
$ws.Cells.Item(68, 11).Value = 2697
$ws.Cells.Item(68, 12).Value = 3157.7142
$ws.Cells.Item(68, 13).Value = -1886
$ws.Cells.Item(68, 14).Value = -4779.7142
$ws.Cells.Item(71, 8).Value = 1033.375
$ws.Cells.Item(71, 9).Value = 899
$ws.Cells.Item(71, 10).Value = 1052.5714
$ws.Cells.Item(71, 11).Value = 8091
$ws.Cells.Item(71, 12).Value = 9473.142600000001
$ws.Cells.Item(71, 13).Value = -4035
$ws.Cells.Item(71, 14).Value = -17585.1426
$ws.Cells.Item(98, 8).Value = 3448.8333
$ws.Cells.Item(98, 9).Value = 0
$ws.Cells.Item(98, 10).Value = 3448.8333
$ws.Cells.Item(98, 11).Value = 0
$ws.Cells.Item(98, 12).Value = 10346.4999
$ws.Cells.Item(98, 13).ClearContents()
$ws.Cells.Item(98, 14).Value = -13342.4999
$ws.Cells.Item(131, 8).Value = 1483.35
$ws.Cells.Item(131, 9).Value = 1390.75
$ws.Cells.Item(131, 10).Value = 1487.2084
$ws.Cells.Item(131, 11).Value = 4172.25
$ws.Cells.Item(131, 12).Value = 4461.6252
$ws.Cells.Item(131, 13).Value = 867.75
$ws.Cells.Item(131, 14).Value = -14541.6252

$ws = $wb.Worksheets("GSM")
$ws.Cells.Item(52, 8).Value = 20605.584
$ws.Cells.Item(52, 9).Value = 19000
$ws.Cells.Item(52, 11).Value = 19000
$ws.Cells.Item(52, 13).Value = -18741
$ws.Cells.Item(80, 8).Value = 2978.75
$ws.Cells.Item(80, 9).Value = 2690
$ws.Cells.Item(80, 11).Value = 2690
$ws.Cells.Item(80, 13).Value = -1692
$ws.Cells.Item(83, 8).Value = 2978.75
$ws.Cells.Item(83, 9).Value = 2690
$ws.Cells.Item(83, 11).Value = 13450
$ws.Cells.Item(83, 13).Value = -8458
$ws.Cells.Item(104, 8).Value = 29890.334
$ws.Cells.Item(104, 10).Value = 29890.334
$ws.Cells.Item(104, 12).Value = 29890.334
$ws.Cells.Item(104, 14).Value = -36878.334
$ws.Cells.Item(126, 8).Value = 7994.4736
$ws.Cells.Item(126, 9).Value = 9780
$ws.Cells.Item(126, 11).Value = 29340
$ws.Cells.Item(126, 13).Value = -26870

$ws = $wb.Worksheets("LTW")
$ws.Cells.Item(46, 8).Value = 1033.3334
$ws.Cells.Item(46, 9).Value = 800
$ws.Cells.Item(46, 10).Value = 1500
$ws.Cells.Item(46, 11).Value = 800
$ws.Cells.Item(46, 12).Value = 1500
$ws.Cells.Item(46, 13).Value = -612
$ws.Cells.Item(46, 14).Value = -1876
$ws.Cells.Item(68, 8).Value = 8077.409
$ws.Cells.Item(68, 10).Value = 10783.111
$ws.Cells.Item(68, 12).Value = 10783.111
$ws.Cells.Item(68, 14).Value = -12281.111
$ws.Cells.Item(71, 8).Value = 8077.409
$ws.Cells.Item(71, 10).Value = 10783.111
$ws.Cells.Item(71, 12).Value = 53915.55500000001
$ws.Cells.Item(71, 14).Value = -61403.55500000001
$ws.Cells.Item(105, 8).Value = 21659.834
$ws.Cells.Item(105, 10).Value = 21659.834
$ws.Cells.Item(105, 12).Value = 21659.834
$ws.Cells.Item(105, 14).Value = -28647.834
$ws.Cells.Item(132, 8).Value = 3680035.2
$ws.Cells.Item(132, 9).Value = 7434.6665
$ws.Cells.Item(132, 11).Value = 22303.9995
$ws.Cells.Item(132, 13).Value = -19773.9995

$ws = $wb.Worksheets("WVR")
$ws.Cells.Item(96, 8).Value = 1374
$ws.Cells.Item(96, 10).Value = 1425.1
$ws.Cells.Item(96, 12).Value = 1425.1
$ws.Cells.Item(96, 14).Value = -4171.1
$ws.Cells.Item(126, 8).Value = 5460.857
$ws.Cells.Item(126, 9).Value = 5650.154
$ws.Cells.Item(126, 10).Value = 3000
$ws.Cells.Item(126, 11).Value = 16950.462
$ws.Cells.Item(126, 12).Value = 9000
$ws.Cells.Item(126, 13).Value = -14480.462
$ws.Cells.Item(126, 14).Value = -13940
